$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Simple value replacements (1-indexed rows, single column table)
$t.Cell(1, 1).Range.Text  = "0M"
$t.Cell(2, 1).Range.Text  = "0M"
$t.Cell(3, 1).Range.Text  = "0M"
$t.Cell(4, 1).Range.Text  = "579"
$t.Cell(6, 1).Range.Text  = "0.01336"
$t.Cell(7, 1).Range.Text  = "0.00228"
$t.Cell(8, 1).Range.Text  = "0.00088"
$t.Cell(9, 1).Range.Text  = "0.00842"
$t.Cell(10, 1).Range.Text = "0.00842"
$t.Cell(11, 1).Range.Text = "0.01336"
$t.Cell(12, 1).Range.Text = "0.13055"

# Collapse the detailed tab-separated breakdown rows down to the single
# summary value that used to live in the earlier rows.
$t.Cell(44, 1).Range.Text = "99.83"
$t.Cell(45, 1).Range.Text = "0.13"
$t.Cell(46, 1).Range.Text = "75"
